# Auto-generated: apply scheduled-runner price/profit refresh to the Anima_Profits workbook.
# For each touched leve row, rewrite the market-price / profit columns (H:N) coming back
# from the price-fetch pass. A few rows lose their M or N cell entirely (no HQ/NQ profit
# could be computed that pass), which are cleared rather than zeroed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2116.8
$ws.Range("I33").Value = 1371
$ws.Range("J33").Value = 5100
$ws.Range("K33").Value = 1371
$ws.Range("L33").Value = 5100
$ws.Range("M33").Value = -1142
$ws.Range("N33").Value = -5558
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H112").Value = 8348.440000000001
$ws.Range("J112").Value = 8348.440000000001
$ws.Range("L112").Value = 25045.32
$ws.Range("N112").Value = -27261.32
$ws.Range("H115").Value = 2489.2856
$ws.Range("I115").Value = 2285
$ws.Range("K115").Value = 6855
$ws.Range("M115").Value = -5288
$ws.Range("H127").Value = 1156.3889
$ws.Range("I127").Value = 671.875
$ws.Range("J127").Value = 1544
$ws.Range("K127").Value = 2015.625
$ws.Range("L127").Value = 4632
$ws.Range("M127").Value = 2944.375
$ws.Range("N127").Value = -14552
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
$ws.Range("H137").Value = 2780.8809
$ws.Range("I137").Value = 2046.3823
$ws.Range("K137").Value = 6139.1469
$ws.Range("M137").Value = -3589.1469
$ws.Range("H138").Value = 2925.8386
$ws.Range("I138").Value = 3322.7693
$ws.Range("J138").Value = 2639.1667
$ws.Range("K138").Value = 9968.3079
$ws.Range("L138").Value = 7917.500100000001
$ws.Range("M138").Value = -4828.3079
$ws.Range("N138").Value = -18197.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4845.364
$ws.Range("I61").Value = 6900
$ws.Range("K61").Value = 6900
$ws.Range("M61").Value = -6688
$ws.Range("H74").Value = 1277.92
$ws.Range("I74").Value = 843.82355
$ws.Range("J74").Value = 2200.375
$ws.Range("K74").Value = 843.82355
$ws.Range("L74").Value = 2200.375
$ws.Range("M74").Value = 30.17645000000005
$ws.Range("N74").Value = -3948.375
$ws.Range("H77").Value = 1277.92
$ws.Range("I77").Value = 843.82355
$ws.Range("J77").Value = 2200.375
$ws.Range("K77").Value = 4219.117749999999
$ws.Range("L77").Value = 11001.875
$ws.Range("M77").Value = 148.8822500000006
$ws.Range("N77").Value = -19737.875
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = $null
$ws.Range("H132").Value = 3097.3076
$ws.Range("I132").Value = 2169.3674
$ws.Range("K132").Value = 6508.1022
$ws.Range("M132").Value = -3978.1022
$ws.Range("H136").Value = 4845.364
$ws.Range("I136").Value = 6900
$ws.Range("K136").Value = 20700
$ws.Range("M136").Value = -18150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2566.641
$ws.Range("I134").Value = 2099.5173
$ws.Range("K134").Value = 6298.5519
$ws.Range("M134").Value = -3763.5519

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7232.3335
$ws.Range("I31").Value = 1317.2609
$ws.Range("J31").Value = 13416.272
$ws.Range("K31").Value = 1317.2609
$ws.Range("L31").Value = 13416.272
$ws.Range("M31").Value = -1022.2609
$ws.Range("N31").Value = -14006.272
$ws.Range("H34").Value = 7232.3335
$ws.Range("I34").Value = 1317.2609
$ws.Range("J34").Value = 13416.272
$ws.Range("K34").Value = 1317.2609
$ws.Range("L34").Value = 13416.272
$ws.Range("M34").Value = -1115.2609
$ws.Range("N34").Value = -13820.272
$ws.Range("H58").Value = 1735.7333
$ws.Range("J58").Value = 2297.6
$ws.Range("L58").Value = 2297.6
$ws.Range("N58").Value = -2703.6
$ws.Range("H132").Value = 6668717.5
$ws.Range("I132").Value = 1839.4706
$ws.Range("J132").Value = 20835834
$ws.Range("K132").Value = 5518.4118
$ws.Range("L132").Value = 62507502
$ws.Range("M132").Value = -2988.4118
$ws.Range("N132").Value = -62512562
$ws.Range("H134").Value = 4889.4443
$ws.Range("I134").Value = 3200
$ws.Range("K134").Value = 9600
$ws.Range("M134").Value = -7065
$ws.Range("H136").Value = 1735.7333
$ws.Range("J136").Value = 2297.6
$ws.Range("L136").Value = 6892.799999999999
$ws.Range("N136").Value = -11992.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 427.25
$ws.Range("J26").Value = 522.7273
$ws.Range("L26").Value = 1568.1819
$ws.Range("N26").Value = -2144.1819
$ws.Range("H86").Value = 1980
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("H89").Value = 1980
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("H122").Value = 7242.467
$ws.Range("J122").Value = 34332.668
$ws.Range("L122").Value = 308994.012
$ws.Range("N122").Value = -313894.012

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1428.6923
$ws.Range("I97").Value = 1021.75
$ws.Range("J97").Value = 2079.8
$ws.Range("K97").Value = 1021.75
$ws.Range("L97").Value = 2079.8
$ws.Range("M97").Value = -525.75
$ws.Range("N97").Value = -3071.8
$ws.Range("H132").Value = 2963.4583
$ws.Range("I132").Value = 2884.1667
$ws.Range("K132").Value = 8652.500100000001
$ws.Range("M132").Value = -6122.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2588.1052
$ws.Range("I132").Value = 1948.6875
$ws.Range("K132").Value = 5846.0625
$ws.Range("M132").Value = -3316.0625
$ws.Range("H136").Value = 10419111
$ws.Range("I136").Value = 2161.5386
$ws.Range("J136").Value = 55559224
$ws.Range("K136").Value = 6484.6158
$ws.Range("L136").Value = 166677672
$ws.Range("M136").Value = -3934.6158
$ws.Range("N136").Value = -166682772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4699.638
$ws.Range("I96").Value = 2934.818
$ws.Range("K96").Value = 2934.818
$ws.Range("M96").Value = -1561.818
$ws.Range("H132").Value = 6668459
$ws.Range("I132").Value = 1440.15
$ws.Range("K132").Value = 4320.450000000001
$ws.Range("M132").Value = -1790.450000000001
$ws.Range("H136").Value = 2374.2173
$ws.Range("I136").Value = 2180.6553
$ws.Range("J136").Value = 2704.4119
$ws.Range("K136").Value = 6541.965899999999
$ws.Range("L136").Value = 8113.2357
$ws.Range("M136").Value = -3991.965899999999
$ws.Range("N136").Value = -13213.2357
